# Update the "Price" (D) and "Volume(1h)" (E) columns for the crypto rows
# that changed, mirroring the latest scrape. Values are written as plain
# text, matching the inline-string cells already in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells such as "1.001" or "30.615.30" look numeric to Excel, so a
# plain assignment would silently convert them to floating point numbers.
# A leading apostrophe forces a literal-text entry (Excel strips the
# apostrophe itself); ClearFormats() at the end removes the "quote prefix"
# flag that gets set along the way so the cells end up unstyled, same as
# in the original workbook.

$ws.Range("D2").Value = "'30.615.30"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "'1.875.60"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'248.14"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4733"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.2908"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "'0.06481"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'22.11"
$ws.Range("E10").Value = "  +4.94%  "
$ws.Range("D11").Value = "'0.07710"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'0.7397"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "'96.36"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "'1.873.68"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'5.163"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "'273.57"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "'30.648.52"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'0.000007513"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "'2.118.44"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D23").Value = "'5.271"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'6.192"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'9.214"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'163.85"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "'18.74"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "'1.512"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'4.275"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'4.100"
$ws.Range("D34").Value = "'0.04796"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").Value = "'1.122"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'0.6960"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.01850"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'2.750"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'6.237"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").Value = "'73.17"
$ws.Range("E41").Value = "  +4.54%  "
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("D43").Value = "'0.4177"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.8352"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'102.34"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'9.377"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'35.45"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'6.996"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'918.52"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D2:D50").ClearFormats()
